$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.621.15"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "3.085.81"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'516.94"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'142.96"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "'0.435"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "'7.30"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "3.617.09"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "'25.80"
$ws.Range("E14").Value = "  -4.86%  "
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "57.705.23"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "3.090.16"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "'6.12"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").Value = "'13.09"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").Value = "'8.20"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").Value = "'337.82"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'0.503"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").Value = "'65.56"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "'0.172"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").Value = "0.0₃0928"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "'6.46"
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("D29").Value = "'7.13"
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "'20.89"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "'1.17"
$ws.Range("E32").Value = "  -5.46%  "
$ws.Range("D33").Value = "'154.32"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "'27.88"
$ws.Range("E34").Value = "  +7.53%  "
$ws.Range("E35").Value = "  -3.32%  "
$ws.Range("D36").Value = "'5.91"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").Value = "'0.0688"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").Value = "3.129.10"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").Value = "'36.93"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.87"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.672"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "2.290.18"
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("D45").Value = "'0.0252"
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Value = "'20.38"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "'0.948"
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").Value = "'5.90"
$ws.Range("E49").Value = "  -4.97%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").Value = "'0.694"
$ws.Range("E51").Value = "  +0.90%  "
